$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username value in A2 (Fernando301 -> Fernando304)
$ws.Range("A2").Value = "Fernando304"

# Move the active selection to A2
$ws.Range("A2").Select()
